# Update generated output figures (column F - "想去人数") on the
# "展览" and "全部类型" worksheets, matching the refreshed data export.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2009
$ws1.Range("F5").Value = 320
$ws1.Range("F6").Value = 563
$ws1.Range("F8").Value = 2050
$ws1.Range("F9").Value = 10400
$ws1.Range("F14").Value = 398
$ws1.Range("F15").Value = 7254
$ws1.Range("F17").Value = 688
$ws1.Range("F18").Value = 149

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2009
$ws4.Range("F5").Value = 320
$ws4.Range("F6").Value = 563
$ws4.Range("F9").Value = 2050
$ws4.Range("F12").Value = 10400
$ws4.Range("F17").Value = 398
$ws4.Range("F18").Value = 7254
$ws4.Range("F20").Value = 688
$ws4.Range("F21").Value = 149
